$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("2025").Range("A2").Value = 62.7
$wb.Worksheets.Item("2030").Range("A2").Value = 214.5
$wb.Worksheets.Item("2040").Range("A2").Value = 390.5000000000001
$wb.Worksheets.Item("2045").Range("A2").Value = 390.5000000000001
$wb.Worksheets.Item("2050").Range("A2").Value = 390.5000000000001
